# Update cryptos list with latest prices and volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.061.50"
$ws.Range("E2").Value = "'  -2.73%  "

$ws.Range("D3").Value = "'1.866.50"
$ws.Range("E3").Value = "'  -2.06%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.09%  "

$ws.Range("D5").Value = "'306.12"
$ws.Range("E5").Value = "'  -2.18%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  -0.08%  "

$ws.Range("D7").Value = "'0.5162"
$ws.Range("E7").Value = "'  -0.98%  "

$ws.Range("D8").Value = "'0.3768"
$ws.Range("E8").Value = "'  -0.58%  "

$ws.Range("D9").Value = "'0.07165"
$ws.Range("E9").Value = "'  -0.96%  "

$ws.Range("D10").Value = "'0.8893"
$ws.Range("E10").Value = "'  -2.45%  "

$ws.Range("D11").Value = "'20.68"
$ws.Range("E11").Value = "'  -2.74%  "

$ws.Range("D12").Value = "'0.07609"

$ws.Range("D13").Value = "'1.881.81"
$ws.Range("E13").Value = "'  -1.45%  "

$ws.Range("D14").Value = "'5.310"
$ws.Range("E14").Value = "'  -2.59%  "

$ws.Range("D15").Value = "'89.70"
$ws.Range("E15").Value = "'  -2.71%  "

$ws.Range("D16").Value = "'1.001"

$ws.Range("D17").Value = "'0.000008485"
$ws.Range("E17").Value = "'  -2.50%  "

$ws.Range("D18").Value = "'14.09"
$ws.Range("E18").Value = "'  -3.13%  "

$ws.Range("D20").Value = "'27.079.06"
$ws.Range("E20").Value = "'  -2.77%  "

$ws.Range("D21").Value = "'5.031"
$ws.Range("E21").Value = "'  -2.39%  "

$ws.Range("D22").Value = "'2.120.72"
$ws.Range("E22").Value = "'  -3.07%  "

$ws.Range("D23").Value = "'10.52"
$ws.Range("E23").Value = "'  -3.17%  "

$ws.Range("D24").Value = "'6.466"
$ws.Range("E24").Value = "'  -2.63%  "

$ws.Range("D25").Value = "'1.838"
$ws.Range("E25").Value = "'  -1.41%  "

$ws.Range("D26").Value = "'147.68"
$ws.Range("E26").Value = "'  -3.81%  "

$ws.Range("D27").Value = "'17.95"
$ws.Range("E27").Value = "'  -2.11%  "

$ws.Range("D28").Value = "'2.097"
$ws.Range("E28").Value = "'  -3.26%  "

$ws.Range("D29").Value = "'112.75"
$ws.Range("E29").Value = "'  -1.78%  "

$ws.Range("D30").Value = "'4.666"
$ws.Range("E30").Value = "'  -4.00%  "

$ws.Range("D31").Value = "'4.680"
$ws.Range("E31").Value = "'  -3.80%  "

$ws.Range("D32").Value = "'0.09147"
$ws.Range("E32").Value = "'  +1.50%  "

$ws.Range("D33").Value = "'0.05118"
$ws.Range("E33").Value = "'  -3.14%  "

$ws.Range("D34").Value = "'3.068"
$ws.Range("E34").Value = "'  -3.53%  "

$ws.Range("D35").Value = "'1.157"
$ws.Range("E35").Value = "'  -6.34%  "

$ws.Range("D36").Value = "'0.7258"
$ws.Range("E36").Value = "'  -6.98%  "

$ws.Range("E37").Value = "'  -2.92%  "

$ws.Range("E38").Value = "'  -0.31%  "

$ws.Range("D39").Value = "'2.498"
$ws.Range("E39").Value = "'  -4.11%  "

$ws.Range("D40").Value = "'1.074"
$ws.Range("E40").Value = "'  -1.69%  "

$ws.Range("D41").Value = "'0.5331"
$ws.Range("E41").Value = "'  -4.40%  "

$ws.Range("D42").Value = "'6.503"

$ws.Range("D43").Value = "'116.29"
$ws.Range("E43").Value = "'  +0.49%  "

$ws.Range("D44").Value = "'8.295"
$ws.Range("E44").Value = "'  -3.01%  "

$ws.Range("E45").Value = "'  -3.34%  "

$ws.Range("D46").Value = "'0.4635"
$ws.Range("E46").Value = "'  -3.69%  "

$ws.Range("D47").Value = "'0.9998"
$ws.Range("E47").Value = "'  -0.11%  "

$ws.Range("D48").Value = "'9.990"
$ws.Range("E48").Value = "'  -4.74%  "

$ws.Range("D49").Value = "'1.572"
$ws.Range("E49").Value = "'  -2.92%  "

$ws.Range("D50").Value = "'36.56"
$ws.Range("E50").Value = "'  -1.21%  "

$ws.Range("D51").Value = "'63.61"
$ws.Range("E51").Value = "'  -5.07%  "
